$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper to set Price (D) and Volume(1h) (E) for a given row
# Price/Volume columns hold text-formatted values (e.g. "98.32", "43.816.11"),
# so force the Text number format before assigning to prevent Excel from
# auto-converting numeric-looking strings into real numbers.
function Set-Row($row, $price, $volume) {
    if ($null -ne $price) {
        $ws.Cells.Item($row, 4).NumberFormat = "@"
        $ws.Cells.Item($row, 4).Value = $price
    }
    if ($null -ne $volume) {
        $ws.Cells.Item($row, 5).NumberFormat = "@"
        $ws.Cells.Item($row, 5).Value = $volume
    }
}

# Row 2 - Bitcoin
Set-Row 2 "43.790.63" "  +0.31%  "
# Row 3 - Ethereum
Set-Row 3 "2.292.90" "  -1.36%  "
# Row 4 - TetherUSD
Set-Row 4 $null "  -0.04%  "
# Row 5 - Solana
Set-Row 5 "98.32" "  +3.54%  "
# Row 6 - BNB
Set-Row 6 "270.10" "  +0.02%  "
# Row 7 - XRP
Set-Row 7 $null "  +0.12%  "
# Row 8 - USDC
Set-Row 8 $null "  -0.06%  "
# Row 9 - Cardano
Set-Row 9 $null "  -2.08%  "
# Row 10 - Avalanche
Set-Row 10 "45.42" "  +0.46%  "
# Row 11 - Dogecoin
Set-Row 11 "0.0938" "  -0.90%  "
# Row 12 - Polkadot
Set-Row 12 "7.90" "  -2.56%  "
# Row 13 - TRON
Set-Row 13 "0.107" "  +1.47%  "
# Row 14 - Chainlink
Set-Row 14 "15.84" "  +0.88%  "
# Row 15 - WrappedliquidstakedEther2.0
Set-Row 15 "2.636.63" "  -0.95%  "
# Row 16 - Polygon
Set-Row 16 "0.860" "  -0.47%  "
# Row 17 - WrappedEther
Set-Row 17 "2.291.76" "  -1.23%  "
# Row 18 - WrappedBTC
Set-Row 18 "43.791.32" "  +0.51%  "
# Row 19 - ShibaInu
Set-Row 19 $null "  +2.32%  "
# Row 20 - Uniswap
Set-Row 20 $null "  -3.60%  "
# Row 21 - Litecoin
Set-Row 21 "72.43" "  +0.54%  "
# Row 22 - ImmutableX
Set-Row 22 $null "  +7.96%  "
# Row 23 - BitcoinCash
Set-Row 23 "233.49" "  -3.03%  "
# Row 24 - InternetComputer(DFINITY)
Set-Row 24 "9.15" "  -2.92%  "
# Row 25 - PancakeSwap
Set-Row 25 "2.76" "  +9.34%  "
# Row 26 - Dai
Set-Row 26 $null "  -0.01%  "
# Row 27 - Cosmos
Set-Row 27 "11.34" "  -0.89%  "
# Row 28 - WEMIXToken
Set-Row 28 "3.46" "  -0.46%  "
# Row 29 - Toncoin
Set-Row 29 "2.24" "  -2.13%  "
# Row 30 - InjectiveProtocol
Set-Row 30 "38.23" "  -0.22%  "
# Row 31 - Monero
Set-Row 31 "176.54" "  +2.17%  "
# Row 32 - EthereumClassic
Set-Row 32 "21.87" "  -3.24%  "
# Row 33 - Hedera
Set-Row 33 $null "  -1.13%  "
# Row 34 - Filecoin
Set-Row 34 "5.48" "  +0.00%  "
# Row 35 - Stellar
Set-Row 35 $null "  +0.69%  "
# Row 36 - RenderToken
Set-Row 36 "4.69" "  +6.65%  "
# Row 37 - Kaspa
Set-Row 37 $null "  +1.55%  "
# Row 38 - VeChain
Set-Row 38 $null "  -2.22%  "
# Row 39 - NEARProtocol
Set-Row 39 $null "  +4.37%  "

# Rows 40/41 - swap Algorand and LidoDAOToken entries
$ws.Cells.Item(40, 2).Value = "LidoDAOToken"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "2.34"
$ws.Cells.Item(40, 5).NumberFormat = "@"
$ws.Cells.Item(40, 5).Value = "  -0.07%  "

$ws.Cells.Item(41, 2).Value = "Algorand"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "0.237"
$ws.Cells.Item(41, 5).NumberFormat = "@"
$ws.Cells.Item(41, 5).Value = "  +0.16%  "

# Row 42 - ARBITRUM
Set-Row 42 $null "  +1.24%  "
# Row 43 - Celestia
Set-Row 43 "12.19" "  +0.95%  "
# Row 44 - MultiversX
Set-Row 44 "64.64" "  +3.92%  "
# Row 45 - FraxShare
Set-Row 45 "8.83" "  -3.95%  "
# Row 46 - THORChain
Set-Row 46 $null "  -2.19%  "
# Row 47 - Cronos
Set-Row 47 $null "  -0.64%  "
# Row 48 - TrustWalletToken
Set-Row 48 $null "  +1.36%  "
# Row 49 - Aave
Set-Row 49 "98.82" "  -1.86%  "
# Row 50 - WOONetwork
Set-Row 50 "0.440" "  +5.63%  "
# Row 51 - Stacks
Set-Row 51 $null "  +11.12%  "
